# Minor alterations of heatmaps and italicizing gene names
#
# 1. Lower-case the four T6SS gene-name labels in column G of the "T6SS"
#    sheet (TssA1/TssC1/Hcp1/ClpV1 -> tssA1/tssC1/hcp1/clpV1).
# 2. Switch the workbook's active tab from "Phage Genes" to "T6SS" and
#    move each sheet's selected cell to match.

$wb = $excel.ActiveWorkbook

# --- 1. Gene name relabeling on the T6SS sheet -----------------------------
$t6ss = $wb.Worksheets.Item("T6SS")
$t6ss.Range("G2").Value = "tssA1"
$t6ss.Range("G3").Value = "tssC1"
$t6ss.Range("G4").Value = "hcp1"
$t6ss.Range("G5").Value = "clpV1"

# --- 2. Active sheet / selection bookkeeping -------------------------------
# "Phage Genes" was the active tab; it loses that status and its selection
# moves to N8.
$phage = $wb.Worksheets.Item("Phage Genes")
$phage.Range("N8").Select() | Out-Null

# "T6SS" becomes the active tab, with G5 selected.
$t6ss.Activate() | Out-Null
$t6ss.Range("G5").Select() | Out-Null
